$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old row 6 entries (27 / 6 / "started to hate slides :)") ---
# These values move elsewhere (row 10) in a different shape, so just clear this row.
$ws.Range("A6:D6").ClearContents()

# --- Plain numeric "days-ago" column cells that keep the same row# but change value ---
$ws.Range("A8").Value = 28
$ws.Range("A9").Value = 28
$ws.Range("A10").Value = 27
$ws.Range("A11").Value = 26
$ws.Range("A12").Value = 25
$ws.Range("A13").Value = 24
$ws.Range("A14").Value = 24
$ws.Range("A15").Value = 22
$ws.Range("A16").Value = 21

$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 3
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 3
$ws.Range("B13").Value = 3
$ws.Range("B14").Value = 3
$ws.Range("B15").Value = 3
$ws.Range("B16").Value = 3

# --- Summary rows at the bottom (text labels re-use existing shared strings) ---
$ws.Range("A18").Value = "before"
$ws.Range("B18").Value = 18
$ws.Range("A19").Value = "writing proposal "
$ws.Range("B19").Value = 18
$ws.Range("A20").Value = "Unit"
$ws.Range("B20").Value = 8
$ws.Range("A24").Value = "sebastien"
$ws.Range("B24").Value = 250

# --- New hour entries ---
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 2

# --- Numeric-looking text entries "6.5" / "1.5" must land as shared-string text,   ---
# --- not be auto-coerced to numbers, and must NOT leave a stray cell style behind. ---
# --- Stage them through A4 (temporarily Text-formatted), copy/paste VALUES ONLY    ---
# --- into B10/B9, which carries the string without carrying the temp formatting.   ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "6.5"
$ws.Range("A4").Copy()
$ws.Range("B10").PasteSpecial(-4104)
$ws.Range("D10").Value = "started to hate slides :)"

$ws.Range("A4").Value = "1.5"
$ws.Range("A4").Copy()
$ws.Range("B9").PasteSpecial(-4104)

# --- New free-text notes (fresh shared strings, inserted in this order) ---
$ws.Range("C4").Value = "stef welcome"
$ws.Range("C5").Value = "discussion with damien on the phone"

# --- Now repurpose A4 (still holding the temp Text format) for its real content: ---
# --- the two date entries, which reuse this exact same cell-style slot.          ---
$ws.Range("A4:A5").NumberFormat = "d-mmm"
$ws.Range("A4").Value = 42192
$ws.Range("A5").Value = 42192
$ws.Range("A7").NumberFormat = "d-mmm"
$ws.Range("A7").Value = 42185

# --- Totals formula now covers the extended data range ---
$ws.Range("C1").Formula = "=SUM(B6:B29)"

# --- View bits ---
$ws.Application.ActiveWindow.RangeSelection
$ws.Range("A8").Select()
